# Delete row 589 ("「デザイン」..." post) entirely, causing all
# subsequent rows to shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(589).Delete()
